$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing company name string (B1) and the tree name in B27
$ws.Range("B1").Value = "ООО Ромашка"
$ws.Range("B5").Value = "Береза"
$ws.Range("B27").Value = "Дуб"

# Move the active selection to B26 as recorded in the saved sheet view
$ws.Range("B26").Select()
